# Insert a new weekly record row above row 54, shifting existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54; rows 54..175 shift down to 55..176.
$ws.Rows(54).Insert()

# Fill the new row 54 with the new weekly record.
$ws.Range("A54").Value = 5
$ws.Range("B54").Value = "Macroferia Regional de Talca"
$ws.Range("C54").Value = "Maule"
$ws.Range("D54").Value = 44526
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = 100112008
$ws.Range("G54").Value = "Coliflor"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 4000
$ws.Range("K54").Value = 600
$ws.Range("L54").Value = 600
$ws.Range("M54").Value = 600
$ws.Range("N54").Value = "`$/unidad"
$ws.Range("O54").Value = "Región del Maule"
$ws.Range("P54").Value = 600
$ws.Range("Q54").Value = 1
$ws.Range("R54").Value = "Hortaliza"
